# Insert a new weekly price record for "Acelga" (Feria Lagunitas de Puerto
# Montt) as row 93, pushing the existing rows 93-110 down to 94-111.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 93..110 down to 94..111 (new blank row appears at 93).
$ws.Rows.Item(93).Insert()

# Populate the newly inserted row 93 with the new observation.
$ws.Cells.Item(93, 1).Value  = 4
$ws.Cells.Item(93, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(93, 3).Value  = "Los Lagos"
$ws.Cells.Item(93, 4).Value  = 44474
$ws.Cells.Item(93, 5).Value  = 10
$ws.Cells.Item(93, 6).Value  = 100112009
$ws.Cells.Item(93, 7).Value  = "Acelga"
$ws.Cells.Item(93, 8).Value  = "Sin especificar"
$ws.Cells.Item(93, 9).Value  = "Primera"
$ws.Cells.Item(93, 10).Value = 200
$ws.Cells.Item(93, 11).Value = 4000
$ws.Cells.Item(93, 12).Value = 4000
$ws.Cells.Item(93, 13).Value = 4000
$ws.Cells.Item(93, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(93, 15).Value = "Región del Maule"
$ws.Cells.Item(93, 16).Value = 1000
$ws.Cells.Item(93, 17).Value = 4
$ws.Cells.Item(93, 18).Value = "Hortaliza"
